# Commit: regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# The upstream pipeline recomputed column G ("K") for each row of the
# season log using the new K-based definition instead of the old Strike#-based
# one. This updates the recomputed K values in place, row by row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new K (column G) value, taken from the regenerated
# save_data output.
$kUpdates = @{
    2 = 1
    3 = 0
    4 = 0
    5 = 1
    6 = 2
    7 = 2
    8 = 2
    9 = 2
    10 = 1
    11 = 1
    12 = 1
    13 = 2
    14 = 1
    15 = 0
    16 = 1
    18 = 4
    19 = 3
    20 = 3
    21 = 1
    23 = 1
    24 = 1
    25 = 2
    26 = 0
    27 = 1
    28 = 1
    29 = 0
    30 = 0
    31 = 0
    32 = 2
    33 = 1
    34 = 2
    35 = 1
    36 = 1
    37 = 0
    38 = 0
    39 = 0
    40 = 0
    41 = 1
    42 = 1
    43 = 2
    44 = 1
    45 = 3
    46 = 3
    47 = 1
    48 = 0
    49 = 0
    50 = 2
    51 = 1
    52 = 2
    53 = 1
    54 = 0
    55 = 1
    56 = 1
    58 = 2
    59 = 2
    60 = 2
    61 = 3
    62 = 0
    63 = 1
    64 = 0
    65 = 2
    66 = 1
    68 = 3
    69 = 2
    72 = 1
    73 = 1
}

foreach ($row in $kUpdates.Keys) {
    $ws.Cells.Item($row, 7).Value = $kUpdates[$row]
}
